$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1251.6492
$ws.Range("I15").Value = 1251.6492
$ws.Range("K15").Value = 3754.9476
$ws.Range("M15").Value = -3585.9476
$ws.Range("H16").Value = 509
$ws.Range("I16").Value = 509
$ws.Range("K16").Value = 509
$ws.Range("M16").Value = -279
$ws.Range("H33").Value = 236350.16
$ws.Range("I33").Value = 129.14285
$ws.Range("K33").Value = 129.14285
$ws.Range("M33").Value = 99.85714999999999
$ws.Range("H137").Value = 1800.5186
$ws.Range("I137").Value = 1870.2667
$ws.Range("J137").Value = 1713.3334
$ws.Range("K137").Value = 5610.800099999999
$ws.Range("L137").Value = 5140.0002
$ws.Range("M137").Value = -3060.800099999999
$ws.Range("N137").Value = -10240.0002
$ws.Range("H138").Value = 2909
$ws.Range("J138").Value = 3129.3635
$ws.Range("L138").Value = 9388.0905
$ws.Range("N138").Value = -19668.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 72710.78999999999
$ws.Range("I2").Value = 1412.4546
$ws.Range("J2").Value = 334138
$ws.Range("K2").Value = 1412.4546
$ws.Range("L2").Value = 334138
$ws.Range("M2").Value = -1299.4546
$ws.Range("N2").Value = -334364
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("K30").Value = 1000
$ws.Range("M30").Value = -850
$ws.Range("H32").Value = 34684.98
$ws.Range("I32").Value = 7438.479
$ws.Range("J32").Value = 221518.14
$ws.Range("K32").Value = 7438.479
$ws.Range("L32").Value = 221518.14
$ws.Range("M32").Value = -7151.479
$ws.Range("N32").Value = -222092.14
$ws.Range("H59").Value = 16000
$ws.Range("J59").Value = 16000
$ws.Range("L59").Value = 16000
$ws.Range("N59").Value = -17608
$ws.Range("H116").Value = 72710.78999999999
$ws.Range("I116").Value = 1412.4546
$ws.Range("J116").Value = 334138
$ws.Range("K116").Value = 1412.4546
$ws.Range("L116").Value = 334138
$ws.Range("M116").Value = 881.5454
$ws.Range("N116").Value = -338726
$ws.Range("H117").Value = 33293.332
$ws.Range("J117").Value = 33293.332
$ws.Range("L117").Value = 33293.332
$ws.Range("N117").Value = -42471.332
$ws.Range("H135").Value = 42884
$ws.Range("J135").Value = 42884
$ws.Range("L135").Value = 42884
$ws.Range("N135").Value = -53024

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 72710.78999999999
$ws.Range("I3").Value = 1412.4546
$ws.Range("J3").Value = 334138
$ws.Range("K3").Value = 1412.4546
$ws.Range("L3").Value = 334138
$ws.Range("M3").Value = -1298.4546
$ws.Range("N3").Value = -334366
$ws.Range("H86").Value = 55348.57
$ws.Range("I86").Value = 103791.63
$ws.Range("K86").Value = 103791.63
$ws.Range("M86").Value = -102668.63
$ws.Range("H89").Value = 55348.57
$ws.Range("I89").Value = 103791.63
$ws.Range("K89").Value = 518958.15
$ws.Range("M89").Value = -513342.15
$ws.Range("H99").Value = 2356.3572
$ws.Range("I99").Value = 2222.25
$ws.Range("J99").Value = 2410
$ws.Range("K99").Value = 2222.25
$ws.Range("L99").Value = 2410
$ws.Range("M99").Value = -724.25
$ws.Range("N99").Value = -5406
$ws.Range("H134").Value = 2716.524
$ws.Range("J134").Value = 1753.5
$ws.Range("L134").Value = 5260.5
$ws.Range("N134").Value = -10330.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 9474.5
$ws.Range("I45").Value = 7632.6665
$ws.Range("K45").Value = 7632.6665
$ws.Range("M45").Value = -7039.6665
$ws.Range("H58").Value = 1504.6666
$ws.Range("I58").Value = 1100
$ws.Range("J58").Value = 1707
$ws.Range("K58").Value = 1100
$ws.Range("L58").Value = 1707
$ws.Range("M58").Value = -897
$ws.Range("N58").Value = -2113
$ws.Range("H132").Value = 150006850
$ws.Range("J132").Value = 71433040
$ws.Range("L132").Value = 214299120
$ws.Range("N132").Value = -214304180
$ws.Range("H136").Value = 1504.6666
$ws.Range("I136").Value = 1100
$ws.Range("J136").Value = 1707
$ws.Range("K136").Value = 3300
$ws.Range("L136").Value = 5121
$ws.Range("M136").Value = -750
$ws.Range("N136").Value = -10221

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 3682.8572
$ws.Range("I56").Value = 3682.8572
$ws.Range("K56").Value = 3682.8572
$ws.Range("M56").Value = -3152.8572
$ws.Range("H117").Value = 6853.8945
$ws.Range("J117").Value = 7558.4116
$ws.Range("L117").Value = 22675.2348
$ws.Range("N117").Value = -29559.2348
$ws.Range("H131").Value = 836.05
$ws.Range("J131").Value = 850.05206
$ws.Range("L131").Value = 2550.15618
$ws.Range("N131").Value = -12630.15618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 11993.333
$ws.Range("I21").Value = 11980
$ws.Range("J21").Value = 12000
$ws.Range("K21").Value = 11980
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = -11807
$ws.Range("N21").Value = -12346
$ws.Range("H30").Value = 11993.333
$ws.Range("I30").Value = 11980
$ws.Range("J30").Value = 12000
$ws.Range("K30").Value = 11980
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = -11875
$ws.Range("N30").Value = -12210
$ws.Range("H95").Value = 19900
$ws.Range("J95").Value = 19900
$ws.Range("L95").Value = 19900
$ws.Range("N95").Value = -25392
$ws.Range("H97").Value = 142859620
$ws.Range("I97").Value = 142859620
$ws.Range("K97").Value = 142859620
$ws.Range("M97").Value = -142859124

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1359.1724
$ws.Range("J22").Value = 973.5909
$ws.Range("L22").Value = 973.5909
$ws.Range("N22").Value = -1563.5909
$ws.Range("H27").Value = 1359.1724
$ws.Range("J27").Value = 973.5909
$ws.Range("L27").Value = 973.5909
$ws.Range("N27").Value = -1187.5909
$ws.Range("H57").Value = 47046
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 47046
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 47046
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -48178
$ws.Range("H82").Value = 2591.4285
$ws.Range("I82").Value = 2720
$ws.Range("J82").Value = 2495
$ws.Range("K82").Value = 2720
$ws.Range("L82").Value = 2495
$ws.Range("M82").Value = -2359
$ws.Range("N82").Value = -3217
$ws.Range("H85").Value = 2591.4285
$ws.Range("I85").Value = 2720
$ws.Range("J85").Value = 2495
$ws.Range("K85").Value = 2720
$ws.Range("L85").Value = 2495
$ws.Range("M85").Value = -1472
$ws.Range("N85").Value = -4991
$ws.Range("H100").Value = 2582.7144
$ws.Range("I100").Value = 2199.6667
$ws.Range("J100").Value = 2870
$ws.Range("K100").Value = 2199.6667
$ws.Range("L100").Value = 2870
$ws.Range("M100").Value = -1658.6667
$ws.Range("N100").Value = -3952
$ws.Range("H132").Value = 4604.4546
$ws.Range("I132").Value = 4489.4736
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 13468.4208
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -10938.4208
$ws.Range("N132").Value = -21057.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 31527.428
$ws.Range("J56").Value = 36115.332
$ws.Range("L56").Value = 36115.332
$ws.Range("N56").Value = -37543.332
$ws.Range("H81").Value = 200773.2
$ws.Range("I81").Value = 629
$ws.Range("K81").Value = 1258
$ws.Range("M81").Value = -197
$ws.Range("H84").Value = 200773.2
$ws.Range("I84").Value = 629
$ws.Range("K84").Value = 6290
$ws.Range("M84").Value = -986
